# Refresh the cryptocurrency price / 1h-volume-change table with the
# latest scrape (GitHub Actions scheduled run). Row 33/34 also swap
# (Cosmos overtakes OKB in the ranking), so those two rows get their
# Coin/Link/Price/Volume columns rewritten rather than just the numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Column D/E store plain text (prices/percentages formatted upstream),
    # not numbers. Force text so strings like "1.00" or "0.130" keep their
    # exact digits instead of Excel coercing them into numeric literals,
    # then drop back to the default style so no formatting residue is left
    # behind on the cell (matches the original "General" / unstyled cells).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "66.627.49"
Set-TextValue $ws.Range("E2") "  -3.10%  "
Set-TextValue $ws.Range("D3") "3.573.75"
Set-TextValue $ws.Range("E3") "  -3.40%  "
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "589.27"
Set-TextValue $ws.Range("E5") "  -0.98%  "
Set-TextValue $ws.Range("D6") "183.12"
Set-TextValue $ws.Range("E6") "  +1.00%  "
Set-TextValue $ws.Range("D7") "3.567.65"
Set-TextValue $ws.Range("E7") "  -3.52%  "
Set-TextValue $ws.Range("D8") "0.609"
Set-TextValue $ws.Range("E8") "  -3.33%  "
Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  -0.01%  "
Set-TextValue $ws.Range("D10") "0.669"
Set-TextValue $ws.Range("E10") "  -6.76%  "
Set-TextValue $ws.Range("D11") "53.81"
Set-TextValue $ws.Range("E11") "  -4.39%  "
Set-TextValue $ws.Range("E12") "  -11.08%  "
Set-TextValue $ws.Range("D13") "0.0000251"
Set-TextValue $ws.Range("E13") "  -14.78%  "
Set-TextValue $ws.Range("D14") "9.80"
Set-TextValue $ws.Range("E14") "  -8.46%  "
Set-TextValue $ws.Range("D15") "4.141.81"
Set-TextValue $ws.Range("E15") "  -3.18%  "
Set-TextValue $ws.Range("D16") "3.577.14"
Set-TextValue $ws.Range("E16") "  -3.32%  "
Set-TextValue $ws.Range("D17") "0.126"
Set-TextValue $ws.Range("E17") "  -0.58%  "
Set-TextValue $ws.Range("E18") "  -5.82%  "
Set-TextValue $ws.Range("D19") "66.444.25"
Set-TextValue $ws.Range("E19") "  -2.96%  "
Set-TextValue $ws.Range("D20") "12.15"
Set-TextValue $ws.Range("E20") "  -5.88%  "
Set-TextValue $ws.Range("E21") "  -6.68%  "
Set-TextValue $ws.Range("D22") "392.88"
Set-TextValue $ws.Range("E22") "  -4.95%  "
Set-TextValue $ws.Range("D23") "4.29"
Set-TextValue $ws.Range("E23") "  -6.83%  "
Set-TextValue $ws.Range("E24") "  -5.00%  "
Set-TextValue $ws.Range("E25") "  -5.91%  "
Set-TextValue $ws.Range("D26") "12.30"
Set-TextValue $ws.Range("E26") "  -3.43%  "
Set-TextValue $ws.Range("E27") "  -0.78%  "
Set-TextValue $ws.Range("D28") "10.25"
Set-TextValue $ws.Range("E28") "  -5.81%  "
Set-TextValue $ws.Range("D29") "3.59"
Set-TextValue $ws.Range("E29") "  -7.80%  "
Set-TextValue $ws.Range("D30") "8.92"
Set-TextValue $ws.Range("E30") "  -7.65%  "
Set-TextValue $ws.Range("D31") "31.07"
Set-TextValue $ws.Range("E31") "  -5.88%  "
Set-TextValue $ws.Range("D32") "6.76"
Set-TextValue $ws.Range("E32") "  -8.84%  "
Set-TextValue $ws.Range("B33") "Cosmos"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D33") "11.95"
Set-TextValue $ws.Range("E33") "  -4.38%  "
Set-TextValue $ws.Range("B34") "OKB"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D34") "65.31"
Set-TextValue $ws.Range("E34") "  +0.67%  "
Set-TextValue $ws.Range("D35") "611.01"
Set-TextValue $ws.Range("E35") "  +1.59%  "
Set-TextValue $ws.Range("E36") "  -6.49%  "
Set-TextValue $ws.Range("D37") "41.42"
Set-TextValue $ws.Range("E37") "  -5.05%  "
Set-TextValue $ws.Range("E39") "  -0.12%  "
Set-TextValue $ws.Range("E40") "  -8.02%  "
Set-TextValue $ws.Range("D41") "0.0₃0742"
Set-TextValue $ws.Range("E41") "  -16.15%  "
Set-TextValue $ws.Range("D42") "0.130"
Set-TextValue $ws.Range("E42") "  -6.18%  "
Set-TextValue $ws.Range("D43") "2.922.30"
Set-TextValue $ws.Range("E43") "  +6.71%  "
Set-TextValue $ws.Range("E44") "  -9.25%  "
Set-TextValue $ws.Range("E45") "  -8.38%  "
Set-TextValue $ws.Range("D46") "2.40"
Set-TextValue $ws.Range("E46") "  -10.14%  "
Set-TextValue $ws.Range("E47") "  -4.17%  "
Set-TextValue $ws.Range("D48") "3.05"
Set-TextValue $ws.Range("E48") "  -0.96%  "
Set-TextValue $ws.Range("D49") "136.93"
Set-TextValue $ws.Range("E49") "  -2.62%  "
Set-TextValue $ws.Range("E50") "  -8.26%  "
Set-TextValue $ws.Range("D51") "8.24"
Set-TextValue $ws.Range("E51") "  -9.60%  "
